$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Mapping of row -> new DAMSLTag (col I) and new DialogAct (col J)
$updates = @(
    @{ Row = 3; DamslTag = "b"; DialogAct = "Acknowledge (Backchannel)" },
    @{ Row = 10; DamslTag = "sv"; DialogAct = "Statement-opinion" },
    @{ Row = 20; DamslTag = "sd"; DialogAct = "Statement-non-opinion" },
    @{ Row = 23; DamslTag = "%"; DialogAct = "Uninterpretable" },
    @{ Row = 27; DamslTag = "sv"; DialogAct = "Statement-opinion" },
    @{ Row = 30; DamslTag = "sv"; DialogAct = "Statement-opinion" },
    @{ Row = 35; DamslTag = "b"; DialogAct = "Acknowledge (Backchannel)" },
    @{ Row = 39; DamslTag = "b"; DialogAct = "Acknowledge (Backchannel)" },
    @{ Row = 45; DamslTag = "sd"; DialogAct = "Statement-non-opinion" },
    @{ Row = 49; DamslTag = "b"; DialogAct = "Acknowledge (Backchannel)" },
    @{ Row = 57; DamslTag = "sd"; DialogAct = "Statement-non-opinion" },
    @{ Row = 68; DamslTag = "sv"; DialogAct = "Statement-opinion" },
    @{ Row = 88; DamslTag = "aa"; DialogAct = "Agree/Accept" },
    @{ Row = 90; DamslTag = "aa"; DialogAct = "Agree/Accept" },
    @{ Row = 103; DamslTag = "sd"; DialogAct = "Statement-non-opinion" },
    @{ Row = 115; DamslTag = "b"; DialogAct = "Acknowledge (Backchannel)" },
    @{ Row = 117; DamslTag = "b"; DialogAct = "Acknowledge (Backchannel)" },
    @{ Row = 120; DamslTag = "sd"; DialogAct = "Statement-non-opinion" },
    @{ Row = 121; DamslTag = "b"; DialogAct = "Acknowledge (Backchannel)" },
    @{ Row = 128; DamslTag = "aa"; DialogAct = "Agree/Accept" },
    @{ Row = 129; DamslTag = "aa"; DialogAct = "Agree/Accept" },
    @{ Row = 160; DamslTag = "b"; DialogAct = "Acknowledge (Backchannel)" },
    @{ Row = 165; DamslTag = "sd"; DialogAct = "Statement-non-opinion" },
    @{ Row = 175; DamslTag = "sv"; DialogAct = "Statement-opinion" },
    @{ Row = 184; DamslTag = "sv"; DialogAct = "Statement-opinion" },
    @{ Row = 185; DamslTag = "aa"; DialogAct = "Agree/Accept" },
    @{ Row = 186; DamslTag = "aa"; DialogAct = "Agree/Accept" },
    @{ Row = 187; DamslTag = "sd"; DialogAct = "Statement-non-opinion" },
    @{ Row = 192; DamslTag = "sd"; DialogAct = "Statement-non-opinion" },
    @{ Row = 193; DamslTag = "sd"; DialogAct = "Statement-non-opinion" },
    @{ Row = 195; DamslTag = "aa"; DialogAct = "Agree/Accept" },
    @{ Row = 202; DamslTag = "sd"; DialogAct = "Statement-non-opinion" },
    @{ Row = 206; DamslTag = "sv"; DialogAct = "Statement-opinion" },
    @{ Row = 208; DamslTag = "sd"; DialogAct = "Statement-non-opinion" },
    @{ Row = 226; DamslTag = "sd"; DialogAct = "Statement-non-opinion" },
    @{ Row = 256; DamslTag = "sv"; DialogAct = "Statement-opinion" },
    @{ Row = 262; DamslTag = "%"; DialogAct = "Uninterpretable" },
    @{ Row = 277; DamslTag = "b"; DialogAct = "Acknowledge (Backchannel)" },
    @{ Row = 280; DamslTag = "sd"; DialogAct = "Statement-non-opinion" },
    @{ Row = 281; DamslTag = "sd"; DialogAct = "Statement-non-opinion" },
    @{ Row = 282; DamslTag = "sd"; DialogAct = "Statement-non-opinion" },
    @{ Row = 316; DamslTag = "sv"; DialogAct = "Statement-opinion" },
    @{ Row = 325; DamslTag = "sv"; DialogAct = "Statement-opinion" },
    @{ Row = 330; DamslTag = "sv"; DialogAct = "Statement-opinion" },
    @{ Row = 333; DamslTag = "b"; DialogAct = "Acknowledge (Backchannel)" },
    @{ Row = 339; DamslTag = "sd"; DialogAct = "Statement-non-opinion" },
    @{ Row = 343; DamslTag = "aa"; DialogAct = "Agree/Accept" },
    @{ Row = 355; DamslTag = "sd"; DialogAct = "Statement-non-opinion" },
    @{ Row = 362; DamslTag = "aa"; DialogAct = "Agree/Accept" },
    @{ Row = 371; DamslTag = "%"; DialogAct = "Uninterpretable" },
    @{ Row = 372; DamslTag = "sd"; DialogAct = "Statement-non-opinion" },
    @{ Row = 374; DamslTag = "sd"; DialogAct = "Statement-non-opinion" },
    @{ Row = 375; DamslTag = "sd"; DialogAct = "Statement-non-opinion" },
    @{ Row = 378; DamslTag = "sv"; DialogAct = "Statement-opinion" },
    @{ Row = 406; DamslTag = "%"; DialogAct = "Uninterpretable" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.DamslTag
    $ws.Cells.Item($u.Row, 10).Value = $u.DialogAct
}
